$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/number formats/borders) from column Q (year 2019)
# into the new column R (year 2020), matching each row's existing look.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new 2020 values for column R
$ws.Range("R4").Value = 2020
$ws.Range("R5").Value = 2.1
$ws.Range("R6").Value = 2.4
$ws.Range("R7").Value = 1.4
$ws.Range("R8").Value = 3.2
$ws.Range("R9").Value = 2.4
$ws.Range("R10").Value = 0.8
$ws.Range("R11").Value = 2.2000000000000002
$ws.Range("R12").Value = 4.5
$ws.Range("R13").Value = 1.4
$ws.Range("R14").Value = 3.2

# Update the selection to match the target state
[void]$ws.Range("R16:R17").Select()
